# UP _ exo GAB rename écran code incorrect ressaisie
#
# Slide 5 ("maquettage_GAB_retirer de l'argent dont s'identifier.pptx")
# contains a text box named "ZoneTexte 21" holding the screen caption
# "Ecran : Code confidentiel incorrect". The author renamed it to mention
# the re-entry ("ressaisie") step and widened the (auto-fit) text box to
# fit the longer caption.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item("ZoneTexte 21")

# Widen the textbox to accommodate the longer caption (height is untouched,
# only the width grows from 3624710 EMU to 4676280 EMU; COM works in points,
# 1 pt = 12700 EMU).
$sh.Width = 4676280 / 12700

# Update the caption text itself.
$sh.TextFrame.TextRange.Text = "Ecran : Code confidentiel incorrect : ressaisie"
